# 207; integrate xls with DataProvider
$wb = $excel.ActiveWorkbook

# Work on the "2022" worksheet (sheet2) - add the DataProvider test data
$ws2 = $wb.Worksheets.Item("2022")

# Header row
$ws2.Range("A1").Value = "uid"
$ws2.Range("B1").Value = "pwd"
$ws2.Range("C1").Value = "prodname"

# Data rows
$ws2.Range("A2").Value = "anshika@gmail.com"
$ws2.Range("B2").Value = "Iamking@000"
$ws2.Range("C2").Value = "ADIDAS ORIGINAL"

$ws2.Range("A3").Value = "anshika@gmail.com"
$ws2.Range("B3").Value = "Iamking@000"
$ws2.Range("C3").Value = "ZARA COAT 3"

$ws2.Range("A4").Value = "anshika@gmail.com"
$ws2.Range("B4").Value = "Iamking@000"
$ws2.Range("C4").Value = "IPHONE 13 PRO"

# Apply the new font style (Consolas, 10pt, teal-ish color) to the data rows A2:C4
$dataRange = $ws2.Range("A2:C4")
$dataRange.Font.Name = "Consolas"
$dataRange.Font.Size = 10
$dataRange.Font.Color = 10002730

# Column widths to match sheet1
$ws2.Columns("A:D").ColumnWidth = 14.88671875

# Selection on sheet2
$ws2.Range("C7").Select()

# Make "2022" the active/selected tab
$ws2.Activate()

# Selection on sheet1 changes to full-column A1:XFD1048576 with active cell C21
$ws1 = $wb.Worksheets.Item("2023")
$ws1.Columns.Select() | Out-Null
$ws1.Range("C21").Activate() | Out-Null
